$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 744
$ws.Range("I2").Value = 2157
$ws.Range("J2").Value = 8777
$ws.Range("K2").Value = 51
$ws.Range("L2").Value = 2386
$ws.Range("M2").Value = 134
$ws.Range("N2").Value = 1606
$ws.Range("O2").Value = 8
$ws.Range("Q2").Value = 18
$ws.Range("R2").Value = 119
$ws.Range("S2").Value = 913
$ws.Range("T2").Value = 1573
$ws.Range("U2").Value = 110
$ws.Range("V2").Value = 13252
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 13456
$ws.Range("Y2").Value = 23
$ws.Range("Z2").Value = 214
$ws.Range("AA2").Value = 73
